$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "2024-06-15 01:58:59"
$ws.Range("D8").Value = 200
$ws.Range("E8").Value = 9

# Row 9
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = "2024-06-15 01:59:00"
$ws.Range("D9").Value = 200
$ws.Range("E9").Value = 0
